$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force Age (E) and Phone Number (G) to be stored as text, matching the
# source data (these came in as text, unlike the numeric E2 in the
# existing row).
$ws.Range("E3:E4").NumberFormat = "@"
$ws.Range("G3:G4").NumberFormat = "@"

$ws.Range("A3").Value = "Aarav Mehta"
$ws.Range("B3").Value = "Monday"
$ws.Range("C3").Value = "morning"
$ws.Range("D3").Value = "डॉ. से"
$ws.Range("E3").Value = "28"
$ws.Range("F3").Value = "Male"
$ws.Range("G3").Value = "917823844614"
$ws.Range("H3").Value = "24 MG Road, Bengaluru"
$ws.Range("I3").Value = "2025-06-25 18:43:54"

$ws.Range("A4").Value = "Aarav Mehta"
$ws.Range("B4").Value = "Monday"
$ws.Range("C4").Value = "morning"
$ws.Range("D4").Value = "डॉ. से"
$ws.Range("E4").Value = "28"
$ws.Range("F4").Value = "Male"
$ws.Range("G4").Value = "917823844614"
$ws.Range("H4").Value = "24 MG Road, Bengaluru"
$ws.Range("I4").Value = "2025-06-25 19:32:35"
